$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.887.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.740.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.99%  "

# Row 6
$ws.Range("E6").Value = "  -0.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5171"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.95%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2747"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06159"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.742.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.63%  "

# Row 12
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.598"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "

# Row 16
$ws.Range("E16").Value = "  -0.01%  "

# Row 17
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.903.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006771"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.964.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

# Row 22
$ws.Range("E22").Value = "  +1.97%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.658"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.248"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.508"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.73%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.758"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.71%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.67%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.945"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08288"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.674"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.16%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04596"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.48%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.644"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9871"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6181"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.679"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.924"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.33%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.0000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "

# Row 41
$ws.Range("E41").Value = "  -1.67%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3831"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7407"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.984"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "

# Row 45
$ws.Range("E45").Value = "  +1.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.198"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.90%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05260"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.620"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3401"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
